$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet1: Property/Value metadata updates ---
$ws1.Cells.Item(8, 2).Value = "2025-04-30T13:43:05+00:00"
$ws1.Cells.Item(12, 2).Value = "Model logique d'un lot de soummission"

# --- Sheet2: SubmissionSet structure definition table ---
# Row 2: SubmissionSet
$ws2.Cells.Item(2, 1).Value = "SubmissionSet"
$ws2.Cells.Item(2, 2).Value = "SubmissionSet"
$ws2.Cells.Item(2, 3).Value = $null
$ws2.Cells.Item(2, 4).Value = ""
$ws2.Cells.Item(2, 5).Value = $null
$ws2.Cells.Item(2, 6).Value = "0"
$ws2.Cells.Item(2, 7).Value = "*"
$ws2.Cells.Item(2, 8).Value = ""
$ws2.Cells.Item(2, 9).Value = ""
$ws2.Cells.Item(2, 10).Value = ""
$ws2.Cells.Item(2, 11).Value = "`n"
$ws2.Cells.Item(2, 12).Value = "SubmissionSet (LM)"
$ws2.Cells.Item(2, 13).Value = "Model logique d'un lot de soummission"
$ws2.Cells.Item(2, 14).Value = $null
$ws2.Cells.Item(2, 15).Value = $null
$ws2.Cells.Item(2, 16).Value = ""
$ws2.Cells.Item(2, 17).Value = $null
$ws2.Cells.Item(2, 18).Value = ""
$ws2.Cells.Item(2, 19).Value = ""
$ws2.Cells.Item(2, 20).Value = ""
$ws2.Cells.Item(2, 21).Value = ""
$ws2.Cells.Item(2, 22).Value = ""
$ws2.Cells.Item(2, 23).Value = ""
$ws2.Cells.Item(2, 24).Value = ""
$ws2.Cells.Item(2, 25).Value = ""
$ws2.Cells.Item(2, 26).Value = ""
$ws2.Cells.Item(2, 27).Value = ""
$ws2.Cells.Item(2, 28).Value = ""
$ws2.Cells.Item(2, 29).Value = ""
$ws2.Cells.Item(2, 30).Value = ""
$ws2.Cells.Item(2, 31).Value = ""
$ws2.Cells.Item(2, 32).Value = "Base"
$ws2.Cells.Item(2, 33).Value = "0"
$ws2.Cells.Item(2, 34).Value = "*"
$ws2.Cells.Item(2, 35).Value = ""
$ws2.Cells.Item(2, 36).Value = ""

# Row 3: SubmissionSet.entryUUID
$ws2.Cells.Item(3, 1).Value = "SubmissionSet.entryUUID"
$ws2.Cells.Item(3, 2).Value = "SubmissionSet.entryUUID"
$ws2.Cells.Item(3, 3).Value = $null
$ws2.Cells.Item(3, 4).Value = ""
$ws2.Cells.Item(3, 5).Value = $null
$ws2.Cells.Item(3, 6).Value = "1"
$ws2.Cells.Item(3, 7).Value = "1"
$ws2.Cells.Item(3, 8).Value = ""
$ws2.Cells.Item(3, 9).Value = ""
$ws2.Cells.Item(3, 10).Value = ""
$ws2.Cells.Item(3, 11).Value = "Identifier`n"
$ws2.Cells.Item(3, 12).Value = "Identifiant unique du lot de soumission. Cet attribut est destiné à des fins de gestion interne alors que uniqueId est utilisé à des fins de référence externe. "
$ws2.Cells.Item(3, 13).Value = "Identifiant unique du lot de soumission. Cet attribut est destiné à des fins de gestion interne alors que uniqueId est utilisé à des fins de référence externe."
$ws2.Cells.Item(3, 14).Value = $null
$ws2.Cells.Item(3, 15).Value = $null
$ws2.Cells.Item(3, 16).Value = ""
$ws2.Cells.Item(3, 17).Value = $null
$ws2.Cells.Item(3, 18).Value = ""
$ws2.Cells.Item(3, 19).Value = ""
$ws2.Cells.Item(3, 20).Value = ""
$ws2.Cells.Item(3, 21).Value = ""
$ws2.Cells.Item(3, 22).Value = ""
$ws2.Cells.Item(3, 23).Value = ""
$ws2.Cells.Item(3, 24).Value = ""
$ws2.Cells.Item(3, 25).Value = ""
$ws2.Cells.Item(3, 26).Value = ""
$ws2.Cells.Item(3, 27).Value = ""
$ws2.Cells.Item(3, 28).Value = ""
$ws2.Cells.Item(3, 29).Value = ""
$ws2.Cells.Item(3, 30).Value = ""
$ws2.Cells.Item(3, 31).Value = ""
$ws2.Cells.Item(3, 32).Value = "SubmissionSet.entryUUID"
$ws2.Cells.Item(3, 33).Value = "1"
$ws2.Cells.Item(3, 34).Value = "1"
$ws2.Cells.Item(3, 35).Value = ""
$ws2.Cells.Item(3, 36).Value = ""

# Row 4: SubmissionSet.availabilityStatus
$ws2.Cells.Item(4, 1).Value = "SubmissionSet.availabilityStatus"
$ws2.Cells.Item(4, 2).Value = "SubmissionSet.availabilityStatus"
$ws2.Cells.Item(4, 3).Value = $null
$ws2.Cells.Item(4, 4).Value = ""
$ws2.Cells.Item(4, 5).Value = $null
$ws2.Cells.Item(4, 6).Value = "0"
$ws2.Cells.Item(4, 7).Value = "1"
$ws2.Cells.Item(4, 8).Value = ""
$ws2.Cells.Item(4, 9).Value = ""
$ws2.Cells.Item(4, 10).Value = ""
$ws2.Cells.Item(4, 11).Value = "CodeableConcept`n"
$ws2.Cells.Item(4, 12).Value = "Cette métadonnée représente la pertinence d'un lot de soumission. "
$ws2.Cells.Item(4, 13).Value = "Cette métadonnée représente la pertinence d'un lot de soumission."
$ws2.Cells.Item(4, 14).Value = $null
$ws2.Cells.Item(4, 15).Value = $null
$ws2.Cells.Item(4, 16).Value = ""
$ws2.Cells.Item(4, 17).Value = $null
$ws2.Cells.Item(4, 18).Value = ""
$ws2.Cells.Item(4, 19).Value = ""
$ws2.Cells.Item(4, 20).Value = ""
$ws2.Cells.Item(4, 21).Value = ""
$ws2.Cells.Item(4, 22).Value = ""
$ws2.Cells.Item(4, 23).Value = ""
$ws2.Cells.Item(4, 24).Value = "required"
$ws2.Cells.Item(4, 25).Value = $null
$ws2.Cells.Item(4, 26).Value = "https://mos.esante.gouv.fr/NOS/JDV_J52-AvailabilityStatus-CISIS/FHIR/JDV-J52-AvailabilityStatus-CISIS"
$ws2.Cells.Item(4, 27).Value = ""
$ws2.Cells.Item(4, 28).Value = ""
$ws2.Cells.Item(4, 29).Value = ""
$ws2.Cells.Item(4, 30).Value = ""
$ws2.Cells.Item(4, 31).Value = ""
$ws2.Cells.Item(4, 32).Value = "SubmissionSet.availabilityStatus"
$ws2.Cells.Item(4, 33).Value = "0"
$ws2.Cells.Item(4, 34).Value = "1"
$ws2.Cells.Item(4, 35).Value = ""
$ws2.Cells.Item(4, 36).Value = ""

# Row 5: SubmissionSet.submissionTime
$ws2.Cells.Item(5, 1).Value = "SubmissionSet.submissionTime"
$ws2.Cells.Item(5, 2).Value = "SubmissionSet.submissionTime"
$ws2.Cells.Item(5, 3).Value = $null
$ws2.Cells.Item(5, 4).Value = ""
$ws2.Cells.Item(5, 5).Value = $null
$ws2.Cells.Item(5, 6).Value = "1"
$ws2.Cells.Item(5, 7).Value = "1"
$ws2.Cells.Item(5, 8).Value = ""
$ws2.Cells.Item(5, 9).Value = ""
$ws2.Cells.Item(5, 10).Value = ""
$ws2.Cells.Item(5, 11).Value = "dateTime`n"
$ws2.Cells.Item(5, 12).Value = "Représente la date et heure de soumission."
$ws2.Cells.Item(5, 13).Value = "Représente la date et heure de soumission."
$ws2.Cells.Item(5, 14).Value = $null
$ws2.Cells.Item(5, 15).Value = $null
$ws2.Cells.Item(5, 16).Value = ""
$ws2.Cells.Item(5, 17).Value = $null
$ws2.Cells.Item(5, 18).Value = ""
$ws2.Cells.Item(5, 19).Value = ""
$ws2.Cells.Item(5, 20).Value = ""
$ws2.Cells.Item(5, 21).Value = ""
$ws2.Cells.Item(5, 22).Value = ""
$ws2.Cells.Item(5, 23).Value = ""
$ws2.Cells.Item(5, 24).Value = ""
$ws2.Cells.Item(5, 25).Value = ""
$ws2.Cells.Item(5, 26).Value = ""
$ws2.Cells.Item(5, 27).Value = ""
$ws2.Cells.Item(5, 28).Value = ""
$ws2.Cells.Item(5, 29).Value = ""
$ws2.Cells.Item(5, 30).Value = ""
$ws2.Cells.Item(5, 31).Value = ""
$ws2.Cells.Item(5, 32).Value = "SubmissionSet.submissionTime"
$ws2.Cells.Item(5, 33).Value = "1"
$ws2.Cells.Item(5, 34).Value = "1"
$ws2.Cells.Item(5, 35).Value = ""
$ws2.Cells.Item(5, 36).Value = ""

# Row 6: SubmissionSet.title
$ws2.Cells.Item(6, 1).Value = "SubmissionSet.title"
$ws2.Cells.Item(6, 2).Value = "SubmissionSet.title"
$ws2.Cells.Item(6, 3).Value = $null
$ws2.Cells.Item(6, 4).Value = ""
$ws2.Cells.Item(6, 5).Value = $null
$ws2.Cells.Item(6, 6).Value = "0"
$ws2.Cells.Item(6, 7).Value = "1"
$ws2.Cells.Item(6, 8).Value = ""
$ws2.Cells.Item(6, 9).Value = ""
$ws2.Cells.Item(6, 10).Value = ""
$ws2.Cells.Item(6, 11).Value = "string`n"
$ws2.Cells.Item(6, 12).Value = "Titre du lot de soumission "
$ws2.Cells.Item(6, 13).Value = "Titre du lot de soumission"
$ws2.Cells.Item(6, 14).Value = $null
$ws2.Cells.Item(6, 15).Value = $null
$ws2.Cells.Item(6, 16).Value = ""
$ws2.Cells.Item(6, 17).Value = $null
$ws2.Cells.Item(6, 18).Value = ""
$ws2.Cells.Item(6, 19).Value = ""
$ws2.Cells.Item(6, 20).Value = ""
$ws2.Cells.Item(6, 21).Value = ""
$ws2.Cells.Item(6, 22).Value = ""
$ws2.Cells.Item(6, 23).Value = ""
$ws2.Cells.Item(6, 24).Value = ""
$ws2.Cells.Item(6, 25).Value = ""
$ws2.Cells.Item(6, 26).Value = ""
$ws2.Cells.Item(6, 27).Value = ""
$ws2.Cells.Item(6, 28).Value = ""
$ws2.Cells.Item(6, 29).Value = ""
$ws2.Cells.Item(6, 30).Value = ""
$ws2.Cells.Item(6, 31).Value = ""
$ws2.Cells.Item(6, 32).Value = "SubmissionSet.title"
$ws2.Cells.Item(6, 33).Value = "0"
$ws2.Cells.Item(6, 34).Value = "1"
$ws2.Cells.Item(6, 35).Value = ""
$ws2.Cells.Item(6, 36).Value = ""

# Row 7: SubmissionSet.comments
$ws2.Cells.Item(7, 1).Value = "SubmissionSet.comments"
$ws2.Cells.Item(7, 2).Value = "SubmissionSet.comments"
$ws2.Cells.Item(7, 3).Value = $null
$ws2.Cells.Item(7, 4).Value = ""
$ws2.Cells.Item(7, 5).Value = $null
$ws2.Cells.Item(7, 6).Value = "0"
$ws2.Cells.Item(7, 7).Value = "1"
$ws2.Cells.Item(7, 8).Value = ""
$ws2.Cells.Item(7, 9).Value = ""
$ws2.Cells.Item(7, 10).Value = ""
$ws2.Cells.Item(7, 11).Value = "string`n"
$ws2.Cells.Item(7, 12).Value = "Cette métadonnée contient le commentaire associé au lot de soumission. "
$ws2.Cells.Item(7, 13).Value = "Cette métadonnée contient le commentaire associé au lot de soumission."
$ws2.Cells.Item(7, 14).Value = $null
$ws2.Cells.Item(7, 15).Value = $null
$ws2.Cells.Item(7, 16).Value = ""
$ws2.Cells.Item(7, 17).Value = $null
$ws2.Cells.Item(7, 18).Value = ""
$ws2.Cells.Item(7, 19).Value = ""
$ws2.Cells.Item(7, 20).Value = ""
$ws2.Cells.Item(7, 21).Value = ""
$ws2.Cells.Item(7, 22).Value = ""
$ws2.Cells.Item(7, 23).Value = ""
$ws2.Cells.Item(7, 24).Value = ""
$ws2.Cells.Item(7, 25).Value = ""
$ws2.Cells.Item(7, 26).Value = ""
$ws2.Cells.Item(7, 27).Value = ""
$ws2.Cells.Item(7, 28).Value = ""
$ws2.Cells.Item(7, 29).Value = ""
$ws2.Cells.Item(7, 30).Value = ""
$ws2.Cells.Item(7, 31).Value = ""
$ws2.Cells.Item(7, 32).Value = "SubmissionSet.comments"
$ws2.Cells.Item(7, 33).Value = "0"
$ws2.Cells.Item(7, 34).Value = "1"
$ws2.Cells.Item(7, 35).Value = ""
$ws2.Cells.Item(7, 36).Value = ""

# Row 8: SubmissionSet.patientID
$ws2.Cells.Item(8, 1).Value = "SubmissionSet.patientID"
$ws2.Cells.Item(8, 2).Value = "SubmissionSet.patientID"
$ws2.Cells.Item(8, 3).Value = $null
$ws2.Cells.Item(8, 4).Value = ""
$ws2.Cells.Item(8, 5).Value = $null
$ws2.Cells.Item(8, 6).Value = "1"
$ws2.Cells.Item(8, 7).Value = "1"
$ws2.Cells.Item(8, 8).Value = ""
$ws2.Cells.Item(8, 9).Value = ""
$ws2.Cells.Item(8, 10).Value = ""
$ws2.Cells.Item(8, 11).Value = "Identifier`n"
$ws2.Cells.Item(8, 12).Value = "Cette métadonnée représente l’identifiant du patient, en l’occurrence, le matricule INS (NIR ou NIA) du patient tel que défini dans le cadre juridique. "
$ws2.Cells.Item(8, 13).Value = "Cette métadonnée représente l’identifiant du patient, en l’occurrence, le matricule INS (NIR ou NIA) du patient tel que défini dans le cadre juridique."
$ws2.Cells.Item(8, 14).Value = $null
$ws2.Cells.Item(8, 15).Value = $null
$ws2.Cells.Item(8, 16).Value = ""
$ws2.Cells.Item(8, 17).Value = $null
$ws2.Cells.Item(8, 18).Value = ""
$ws2.Cells.Item(8, 19).Value = ""
$ws2.Cells.Item(8, 20).Value = ""
$ws2.Cells.Item(8, 21).Value = ""
$ws2.Cells.Item(8, 22).Value = ""
$ws2.Cells.Item(8, 23).Value = ""
$ws2.Cells.Item(8, 24).Value = ""
$ws2.Cells.Item(8, 25).Value = ""
$ws2.Cells.Item(8, 26).Value = ""
$ws2.Cells.Item(8, 27).Value = ""
$ws2.Cells.Item(8, 28).Value = ""
$ws2.Cells.Item(8, 29).Value = ""
$ws2.Cells.Item(8, 30).Value = ""
$ws2.Cells.Item(8, 31).Value = ""
$ws2.Cells.Item(8, 32).Value = "SubmissionSet.patientID"
$ws2.Cells.Item(8, 33).Value = "1"
$ws2.Cells.Item(8, 34).Value = "1"
$ws2.Cells.Item(8, 35).Value = ""
$ws2.Cells.Item(8, 36).Value = ""

# Row 9: SubmissionSet.sourceID
$ws2.Cells.Item(9, 1).Value = "SubmissionSet.sourceID"
$ws2.Cells.Item(9, 2).Value = "SubmissionSet.sourceID"
$ws2.Cells.Item(9, 3).Value = $null
$ws2.Cells.Item(9, 4).Value = ""
$ws2.Cells.Item(9, 5).Value = $null
$ws2.Cells.Item(9, 6).Value = "1"
$ws2.Cells.Item(9, 7).Value = "1"
$ws2.Cells.Item(9, 8).Value = ""
$ws2.Cells.Item(9, 9).Value = ""
$ws2.Cells.Item(9, 10).Value = ""
$ws2.Cells.Item(9, 11).Value = "Identifier`n"
$ws2.Cells.Item(9, 12).Value = "Cette métadonnée représente l’identifiant unique global du système émetteur du lot de soumission. "
$ws2.Cells.Item(9, 13).Value = "Cette métadonnée représente l’identifiant unique global du système émetteur du lot de soumission."
$ws2.Cells.Item(9, 14).Value = $null
$ws2.Cells.Item(9, 15).Value = $null
$ws2.Cells.Item(9, 16).Value = ""
$ws2.Cells.Item(9, 17).Value = $null
$ws2.Cells.Item(9, 18).Value = ""
$ws2.Cells.Item(9, 19).Value = ""
$ws2.Cells.Item(9, 20).Value = ""
$ws2.Cells.Item(9, 21).Value = ""
$ws2.Cells.Item(9, 22).Value = ""
$ws2.Cells.Item(9, 23).Value = ""
$ws2.Cells.Item(9, 24).Value = ""
$ws2.Cells.Item(9, 25).Value = ""
$ws2.Cells.Item(9, 26).Value = ""
$ws2.Cells.Item(9, 27).Value = ""
$ws2.Cells.Item(9, 28).Value = ""
$ws2.Cells.Item(9, 29).Value = ""
$ws2.Cells.Item(9, 30).Value = ""
$ws2.Cells.Item(9, 31).Value = ""
$ws2.Cells.Item(9, 32).Value = "SubmissionSet.sourceID"
$ws2.Cells.Item(9, 33).Value = "1"
$ws2.Cells.Item(9, 34).Value = "1"
$ws2.Cells.Item(9, 35).Value = ""
$ws2.Cells.Item(9, 36).Value = ""

# Row 10: SubmissionSet.uniqueID
$ws2.Cells.Item(10, 1).Value = "SubmissionSet.uniqueID"
$ws2.Cells.Item(10, 2).Value = "SubmissionSet.uniqueID"
$ws2.Cells.Item(10, 3).Value = $null
$ws2.Cells.Item(10, 4).Value = ""
$ws2.Cells.Item(10, 5).Value = $null
$ws2.Cells.Item(10, 6).Value = "1"
$ws2.Cells.Item(10, 7).Value = "1"
$ws2.Cells.Item(10, 8).Value = ""
$ws2.Cells.Item(10, 9).Value = ""
$ws2.Cells.Item(10, 10).Value = ""
$ws2.Cells.Item(10, 11).Value = "Identifier`n"
$ws2.Cells.Item(10, 12).Value = "Identifiant unique global affecté à ce lot de soumission par son créateur. Cet attribut est utilisé à des fins de référence externe alors que entryUUID est destiné à des fins de gestion interne.  "
$ws2.Cells.Item(10, 13).Value = "Identifiant unique global affecté à ce lot de soumission par son créateur. Cet attribut est utilisé à des fins de référence externe alors que entryUUID est destiné à des fins de gestion interne."
$ws2.Cells.Item(10, 14).Value = $null
$ws2.Cells.Item(10, 15).Value = $null
$ws2.Cells.Item(10, 16).Value = ""
$ws2.Cells.Item(10, 17).Value = $null
$ws2.Cells.Item(10, 18).Value = ""
$ws2.Cells.Item(10, 19).Value = ""
$ws2.Cells.Item(10, 20).Value = ""
$ws2.Cells.Item(10, 21).Value = ""
$ws2.Cells.Item(10, 22).Value = ""
$ws2.Cells.Item(10, 23).Value = ""
$ws2.Cells.Item(10, 24).Value = ""
$ws2.Cells.Item(10, 25).Value = ""
$ws2.Cells.Item(10, 26).Value = ""
$ws2.Cells.Item(10, 27).Value = ""
$ws2.Cells.Item(10, 28).Value = ""
$ws2.Cells.Item(10, 29).Value = ""
$ws2.Cells.Item(10, 30).Value = ""
$ws2.Cells.Item(10, 31).Value = ""
$ws2.Cells.Item(10, 32).Value = "SubmissionSet.uniqueID"
$ws2.Cells.Item(10, 33).Value = "1"
$ws2.Cells.Item(10, 34).Value = "1"
$ws2.Cells.Item(10, 35).Value = ""
$ws2.Cells.Item(10, 36).Value = ""

# Row 11: SubmissionSet.contentTypeCode
$ws2.Cells.Item(11, 1).Value = "SubmissionSet.contentTypeCode"
$ws2.Cells.Item(11, 2).Value = "SubmissionSet.contentTypeCode"
$ws2.Cells.Item(11, 3).Value = $null
$ws2.Cells.Item(11, 4).Value = ""
$ws2.Cells.Item(11, 5).Value = $null
$ws2.Cells.Item(11, 6).Value = "1"
$ws2.Cells.Item(11, 7).Value = "1"
$ws2.Cells.Item(11, 8).Value = ""
$ws2.Cells.Item(11, 9).Value = ""
$ws2.Cells.Item(11, 10).Value = ""
$ws2.Cells.Item(11, 11).Value = "CodeableConcept`n"
$ws2.Cells.Item(11, 12).Value = "Ensemble de métadonnées représentant le type d’activité associé à l’événement clinique ayant abouti à la constitution du lot de soumission. "
$ws2.Cells.Item(11, 13).Value = "**Submission Set**"
$ws2.Cells.Item(11, 14).Value = $null
$ws2.Cells.Item(11, 15).Value = $null
$ws2.Cells.Item(11, 16).Value = ""
$ws2.Cells.Item(11, 17).Value = $null
$ws2.Cells.Item(11, 18).Value = ""
$ws2.Cells.Item(11, 19).Value = ""
$ws2.Cells.Item(11, 20).Value = ""
$ws2.Cells.Item(11, 21).Value = ""
$ws2.Cells.Item(11, 22).Value = ""
$ws2.Cells.Item(11, 23).Value = ""
$ws2.Cells.Item(11, 24).Value = "required"
$ws2.Cells.Item(11, 25).Value = $null
$ws2.Cells.Item(11, 26).Value = "https://mos.esante.gouv.fr/NOS/JDV_J59-ContentTypeCode-DMP/FHIR/JDV-J59-ContentTypeCode-DMP"
$ws2.Cells.Item(11, 27).Value = ""
$ws2.Cells.Item(11, 28).Value = ""
$ws2.Cells.Item(11, 29).Value = ""
$ws2.Cells.Item(11, 30).Value = ""
$ws2.Cells.Item(11, 31).Value = ""
$ws2.Cells.Item(11, 32).Value = "SubmissionSet.contentTypeCode"
$ws2.Cells.Item(11, 33).Value = "1"
$ws2.Cells.Item(11, 34).Value = "1"
$ws2.Cells.Item(11, 35).Value = ""
$ws2.Cells.Item(11, 36).Value = ""

# Row 12: SubmissionSet.author
$ws2.Cells.Item(12, 1).Value = "SubmissionSet.author"
$ws2.Cells.Item(12, 2).Value = "SubmissionSet.author"
$ws2.Cells.Item(12, 3).Value = $null
$ws2.Cells.Item(12, 4).Value = ""
$ws2.Cells.Item(12, 5).Value = $null
$ws2.Cells.Item(12, 6).Value = "1"
$ws2.Cells.Item(12, 7).Value = "1"
$ws2.Cells.Item(12, 8).Value = ""
$ws2.Cells.Item(12, 9).Value = ""
$ws2.Cells.Item(12, 10).Value = ""
$ws2.Cells.Item(12, 11).Value = "Identifier`n"
$ws2.Cells.Item(12, 12).Value = "Représente la personne physique ou morale et/ou le dispositif auteur d’un lot de soumission"
$ws2.Cells.Item(12, 13).Value = "Représente la personne physique ou morale et/ou le dispositif auteur d’un lot de soumission"
$ws2.Cells.Item(12, 14).Value = $null
$ws2.Cells.Item(12, 15).Value = $null
$ws2.Cells.Item(12, 16).Value = ""
$ws2.Cells.Item(12, 17).Value = $null
$ws2.Cells.Item(12, 18).Value = ""
$ws2.Cells.Item(12, 19).Value = ""
$ws2.Cells.Item(12, 20).Value = ""
$ws2.Cells.Item(12, 21).Value = ""
$ws2.Cells.Item(12, 22).Value = ""
$ws2.Cells.Item(12, 23).Value = ""
$ws2.Cells.Item(12, 24).Value = ""
$ws2.Cells.Item(12, 25).Value = ""
$ws2.Cells.Item(12, 26).Value = ""
$ws2.Cells.Item(12, 27).Value = ""
$ws2.Cells.Item(12, 28).Value = ""
$ws2.Cells.Item(12, 29).Value = ""
$ws2.Cells.Item(12, 30).Value = ""
$ws2.Cells.Item(12, 31).Value = ""
$ws2.Cells.Item(12, 32).Value = "SubmissionSet.author"
$ws2.Cells.Item(12, 33).Value = "1"
$ws2.Cells.Item(12, 34).Value = "1"
$ws2.Cells.Item(12, 35).Value = ""
$ws2.Cells.Item(12, 36).Value = ""

# Row 13: SubmissionSet.homeCommunityID
$ws2.Cells.Item(13, 1).Value = "SubmissionSet.homeCommunityID"
$ws2.Cells.Item(13, 2).Value = "SubmissionSet.homeCommunityID"
$ws2.Cells.Item(13, 3).Value = $null
$ws2.Cells.Item(13, 4).Value = ""
$ws2.Cells.Item(13, 5).Value = $null
$ws2.Cells.Item(13, 6).Value = "0"
$ws2.Cells.Item(13, 7).Value = "1"
$ws2.Cells.Item(13, 8).Value = ""
$ws2.Cells.Item(13, 9).Value = ""
$ws2.Cells.Item(13, 10).Value = ""
$ws2.Cells.Item(13, 11).Value = "Identifier`n"
$ws2.Cells.Item(13, 12).Value = "Cette métadonnée correspond à l’identifiant de la communauté représentée par le système cible si celui-ci offre des fonctionnalités de communication avec d’autres communautés telles que présentées dans le profil XCA d’IHE. Elle n’est pas utilisée par les transactions décrites dans ce volet. "
$ws2.Cells.Item(13, 13).Value = "Cette métadonnée correspond à l’identifiant de la communauté représentée par le système cible si celui-ci offre des fonctionnalités de communication avec d’autres communautés telles que présentées dans le profil XCA d’IHE. Elle n’est pas utilisée par les transactions décrites dans ce volet."
$ws2.Cells.Item(13, 14).Value = $null
$ws2.Cells.Item(13, 15).Value = $null
$ws2.Cells.Item(13, 16).Value = ""
$ws2.Cells.Item(13, 17).Value = $null
$ws2.Cells.Item(13, 18).Value = ""
$ws2.Cells.Item(13, 19).Value = ""
$ws2.Cells.Item(13, 20).Value = ""
$ws2.Cells.Item(13, 21).Value = ""
$ws2.Cells.Item(13, 22).Value = ""
$ws2.Cells.Item(13, 23).Value = ""
$ws2.Cells.Item(13, 24).Value = ""
$ws2.Cells.Item(13, 25).Value = ""
$ws2.Cells.Item(13, 26).Value = ""
$ws2.Cells.Item(13, 27).Value = ""
$ws2.Cells.Item(13, 28).Value = ""
$ws2.Cells.Item(13, 29).Value = ""
$ws2.Cells.Item(13, 30).Value = ""
$ws2.Cells.Item(13, 31).Value = ""
$ws2.Cells.Item(13, 32).Value = "SubmissionSet.homeCommunityID"
$ws2.Cells.Item(13, 33).Value = "0"
$ws2.Cells.Item(13, 34).Value = "1"
$ws2.Cells.Item(13, 35).Value = ""
$ws2.Cells.Item(13, 36).Value = ""

# Row 14: SubmissionSet.intendedRecipient
$ws2.Cells.Item(14, 1).Value = "SubmissionSet.intendedRecipient"
$ws2.Cells.Item(14, 2).Value = "SubmissionSet.intendedRecipient"
$ws2.Cells.Item(14, 3).Value = $null
$ws2.Cells.Item(14, 4).Value = ""
$ws2.Cells.Item(14, 5).Value = $null
$ws2.Cells.Item(14, 6).Value = "0"
$ws2.Cells.Item(14, 7).Value = "*"
$ws2.Cells.Item(14, 8).Value = ""
$ws2.Cells.Item(14, 9).Value = ""
$ws2.Cells.Item(14, 10).Value = ""
$ws2.Cells.Item(14, 11).Value = "Identifier`n"
$ws2.Cells.Item(14, 12).Value = "Cette métadonnée représente les destinataires (structure ou professionnel) auxquels lot de soumission est destiné. Elle n’est pas utilisée par les transactions décrites dans ce volet. "
$ws2.Cells.Item(14, 13).Value = "Cette métadonnée représente les destinataires (structure ou professionnel) auxquels lot de soumission est destiné. Elle n’est pas utilisée par les transactions décrites dans ce volet."
$ws2.Cells.Item(14, 14).Value = $null
$ws2.Cells.Item(14, 15).Value = $null
$ws2.Cells.Item(14, 16).Value = ""
$ws2.Cells.Item(14, 17).Value = $null
$ws2.Cells.Item(14, 18).Value = ""
$ws2.Cells.Item(14, 19).Value = ""
$ws2.Cells.Item(14, 20).Value = ""
$ws2.Cells.Item(14, 21).Value = ""
$ws2.Cells.Item(14, 22).Value = ""
$ws2.Cells.Item(14, 23).Value = ""
$ws2.Cells.Item(14, 24).Value = ""
$ws2.Cells.Item(14, 25).Value = ""
$ws2.Cells.Item(14, 26).Value = ""
$ws2.Cells.Item(14, 27).Value = ""
$ws2.Cells.Item(14, 28).Value = ""
$ws2.Cells.Item(14, 29).Value = ""
$ws2.Cells.Item(14, 30).Value = ""
$ws2.Cells.Item(14, 31).Value = ""
$ws2.Cells.Item(14, 32).Value = "SubmissionSet.intendedRecipient"
$ws2.Cells.Item(14, 33).Value = "0"
$ws2.Cells.Item(14, 34).Value = "*"
$ws2.Cells.Item(14, 35).Value = ""
$ws2.Cells.Item(14, 36).Value = ""

# --- Update dimension / autofit columns A, B, and AF (hidden duplicate) ---
$ws2.Columns.Item(1).EntireColumn.AutoFit() | Out-Null
$ws2.Columns.Item(2).EntireColumn.AutoFit() | Out-Null
$ws2.Columns.Item(32).EntireColumn.AutoFit() | Out-Null
